$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Bank and Cash Accounts" ---
$ws1 = $wb.Worksheets.Item("Bank and Cash Accounts")

$sheet1Data = @{
    2 = @(57327.05, 69089.91, 46505.97, 68042.99000000001, 97647.47, 338613.39)
    3 = @(452802.19, 439948.04, 521104.54, 400785.61, 439761.29, 2254401.67)
    4 = @(-11891.18, -20277.56, -16132.08, -18620.19, -13311.51, -80232.52)
    5 = @(-28564.63, -29994.23, -27644.43, -29930.08, -24965.85, -141099.22)
    6 = @(9349.23, 9672.219999999999, 7688.38, 8537.84, 19512.25, 54759.92)
    7 = @(-196582.08, -188143.44, -185434.74, -205682.74, -189920.85, -965763.85)
    8 = @(-18135.67, -16482.58, -15629.56, -19907.67, -12846.96, -83002.44)
    9 = @(81716.10000000001, 76255.21000000001, 76187.89999999999, 68580.99000000001, 110584.34, 413324.54)
    10 = @(346021.01, 340067.57, 406645.98, 271806.75, 426460.18, 1791001.49)
}

foreach ($row in $sheet1Data.Keys) {
    $vals = $sheet1Data[$row]
    for ($col = 2; $col -le 7; $col++) {
        $ws1.Cells.Item($row, $col).Value = $vals[$col - 2]
    }
}

# --- Sheet 2: "Credit Cards" ---
$ws2 = $wb.Worksheets.Item("Credit Cards")

$sheet2Data = @{
    2 = @(-38915.87, -42885.77, -38293.06, -31954.07, -28973.95, -181022.72)
    3 = @(-36353.05, -30457.38, -38819.41, -38946.09, -34430.07, -179006)
    4 = @(-11866.61, -20344.57, -15438.74, -15008.15, -18061.25, -80719.32000000001)
    5 = @(-42160.42, -40308.66, -31167.05, -37792.51, -39080.4, -190509.04)
    6 = @(-57936.9, -44474.84, -54028.1, -73444.08, -49862.93, -279746.85)
    7 = @(187232.85, 178471.22, 177746.36, 197144.9, 170408.6, 911003.9300000001)
}

foreach ($row in $sheet2Data.Keys) {
    $vals = $sheet2Data[$row]
    for ($col = 2; $col -le 7; $col++) {
        $ws2.Cells.Item($row, $col).Value = $vals[$col - 2]
    }
}
